# Blue Dragon translation tracker update:
# Mark the specific rows below as translated: column C -> "SIM", column D -> "UDS".
# (D is a brand-new column for these rows; C previously held "Não".)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,312,313,314,315,316,317,319,325,336,338,344,351,355,358,402,415,416,417,421,425,433,438,444,451,453,454,457,478,485,501,504,523,536,558,559,560,583,585,586,595,618,619,640,641,647,664,693,699,700,701,702,703,704,741,742,752,776,778,782)

foreach ($r in $rows) {
    # Write column D first so the new shared strings are interned in the same
    # order as the source workbook ("UDS" before "SIM").
    $ws.Cells.Item($r, 4).Value = "UDS"
    $ws.Cells.Item($r, 3).Value = "SIM"
}

# Restore the cursor position left by the author after the edit.
$ws.Range("C322").Select() | Out-Null
